$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1596638655462185
$ws.Range("C2").Value = 0.6078431372549019
$ws.Range("J2").Value = 0.008403361344537815
$ws.Range("P2").Value = 0.15406162464986
$ws.Range("S2").Value = 0.07002801120448179
$ws.Range("B3").Value = 0.01327433628318584
$ws.Range("C3").Value = 0.02654867256637168
$ws.Range("J3").Value = 0.03539823008849557
$ws.Range("P3").Value = 0.7831858407079646
$ws.Range("S3").Value = 0.1415929203539823
$ws.Range("J4").Value = 0.08620689655172414
$ws.Range("P4").Value = 0.6379310344827587
$ws.Range("S4").Value = 0.2758620689655172
$ws.Range("B6").Value = 0.06329113924050633
$ws.Range("D6").Value = 0.004219409282700422
$ws.Range("F6").Value = 0.0379746835443038
$ws.Range("J6").Value = 0.29957805907173
$ws.Range("O6").Value = 0.01265822784810127
$ws.Range("Q6").Value = 0.1983122362869198
$ws.Range("R6").Value = 0.05485232067510549
$ws.Range("S6").Value = 0.3291139240506329
$ws.Range("B7").Value = 0.12
$ws.Range("D7").Value = 0.04
$ws.Range("F7").Value = 0.04
$ws.Range("J7").Value = 0.1257142857142857
$ws.Range("O7").Value = 0.01142857142857143
$ws.Range("Q7").Value = 0.1828571428571429
$ws.Range("R7").Value = 0.1142857142857143
$ws.Range("S7").Value = 0.3657142857142857
$ws.Range("B8").Value = 0.1032148900169205
$ws.Range("D8").Value = 0.02199661590524535
$ws.Range("E8").Value = 0.001692047377326565
$ws.Range("F8").Value = 0.05752961082910321
$ws.Range("J8").Value = 0.09983079526226735
$ws.Range("O8").Value = 0.01692047377326565
$ws.Range("Q8").Value = 0.2301184433164128
$ws.Range("R8").Value = 0.07106598984771574
$ws.Range("S8").Value = 0.3976311336717428
$ws.Range("B9").Value = 0.1238938053097345
$ws.Range("D9").Value = 0.01769911504424779
$ws.Range("F9").Value = 0.06194690265486726
$ws.Range("J9").Value = 0.1017699115044248
$ws.Range("O9").Value = 0.02212389380530973
$ws.Range("Q9").Value = 0.2168141592920354
$ws.Range("R9").Value = 0.084070796460177
$ws.Range("S9").Value = 0.3716814159292036
$ws.Range("B10").Value = 0.1137005649717514
$ws.Range("D10").Value = 0.02471751412429379
$ws.Range("E10").Value = 0.002824858757062147
$ws.Range("F10").Value = 0.06285310734463277
$ws.Range("J10").Value = 0.115819209039548
$ws.Range("O10").Value = 0.01271186440677966
$ws.Range("Q10").Value = 0.2535310734463277
$ws.Range("R10").Value = 0.06991525423728813
$ws.Range("S10").Value = 0.3439265536723164
$ws.Range("G11").Value = 0.1152416356877323
$ws.Range("J11").Value = 0.09293680297397769
$ws.Range("K11").Value = 0.171003717472119
$ws.Range("L11").Value = 0.5836431226765799
$ws.Range("S11").Value = 0.03717472118959108
$ws.Range("G12").Value = 0.7628205128205128
$ws.Range("J12").Value = 0.1858974358974359
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.01923076923076923
$ws.Range("S12").Value = 0.02564102564102564
$ws.Range("G13").Value = 0.7906976744186046
$ws.Range("J13").Value = 0.186046511627907
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("G14").Value = 0.4
$ws.Range("J14").Value = 0.6
$ws.Range("F15").Value = 0.00816326530612245
$ws.Range("H15").Value = 0.163265306122449
$ws.Range("I15").Value = 0.06938775510204082
$ws.Range("J15").Value = 0.3959183673469387
$ws.Range("K15").Value = 0.05306122448979592
$ws.Range("M15").Value = 0.0163265306122449
$ws.Range("O15").Value = 0.04897959183673469
$ws.Range("S15").Value = 0.2448979591836735
$ws.Range("F16").Value = 0.01915708812260536
$ws.Range("H16").Value = 0.2030651340996169
$ws.Range("I16").Value = 0.08045977011494253
$ws.Range("J16").Value = 0.3793103448275862
$ws.Range("K16").Value = 0.09195402298850575
$ws.Range("M16").Value = 0.01532567049808429
$ws.Range("O16").Value = 0.06896551724137931
$ws.Range("S16").Value = 0.1417624521072797
$ws.Range("F17").Value = 0.02243589743589744
$ws.Range("H17").Value = 0.2131410256410256
$ws.Range("I17").Value = 0.09615384615384616
$ws.Range("J17").Value = 0.4022435897435898
$ws.Range("K17").Value = 0.0673076923076923
$ws.Range("M17").Value = 0.01923076923076923
$ws.Range("N17").Value = 0.004807692307692308
$ws.Range("O17").Value = 0.05929487179487179
$ws.Range("S17").Value = 0.1153846153846154
$ws.Range("F18").Value = 0.03141361256544502
$ws.Range("H18").Value = 0.1675392670157068
$ws.Range("I18").Value = 0.08900523560209424
$ws.Range("J18").Value = 0.387434554973822
$ws.Range("K18").Value = 0.07329842931937172
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.1465968586387434
$ws.Range("S18").Value = 0.09424083769633508
$ws.Range("F19").Value = 0.02102102102102102
$ws.Range("H19").Value = 0.25
$ws.Range("I19").Value = 0.08558558558558559
$ws.Range("J19").Value = 0.3708708708708709
$ws.Range("K19").Value = 0.09234234234234234
$ws.Range("M19").Value = 0.01876876876876877
$ws.Range("N19").Value = 0.001501501501501501
$ws.Range("O19").Value = 0.06306306306306306
$ws.Range("S19").Value = 0.09684684684684684

Write-Host "Applied team matrix updates from games pulled march 7"
